# Reorders the player data rows (A2:C18) on Sheet1 to match the target
# workbook state. The underlying (Name, Position, Team) triples are
# unchanged; only their row order is different, so we simply rewrite the
# block of rows with the new ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @("Nick Richards", "C", "Charlotte Hornets"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Kelly Olynyk", "C", "Toronto Raptors"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("Julius Randle", "PF", "Minnesota Timberwolves"),
    @("Herbert Jones", "SF,PF", "New Orleans Pelicans"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Kentavious Caldwell-Pope", "SG,SF", "Orlando Magic"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Malik Monk", "SG,SF", "Sacramento Kings"),
    @("Brandon Miller", "SG,SF", "Charlotte Hornets"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
